# summer 24 week 5 updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 1.29
$ws.Range("D2").Value = 1.23

# Row 3
$ws.Range("B3").Value = 1.55
$ws.Range("E3").Value = 1.32

# Row 4
$ws.Range("B4").Value = 1.51
$ws.Range("C4").Value = 1.44

# Row 5 (values swapped between C5 and D5)
$ws.Range("C5").Value = 1.34
$ws.Range("D5").Value = 1.33

# Row 6
$ws.Range("C6").Value = 1.49
